# "thêm số đt tài xế" - add driver phone number row to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The descriptive note in row 9 (merged/wrapped text) no longer needs the
# tall 90pt row height now that a shorter replacement note is used; shrink
# it down to 36pt.
$ws.Rows.Item(9).RowHeight = 36

# Append the new "số điện thoại tài xế" line as row 34 (#12 in the list).
$ws.Range("A34").Value = 12
$ws.Range("B34").Value = "sdttaixe"

# Leave the selection where the user ended up after adding the new row.
[void]$ws.Range("B35").Select()
